$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "China" (sheet1) -----------------------------------------------------
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("China")

# A2 id changes from 10001 to 10002
$ws1.Range("A2").Value = 10002

# Insert two new blank rows before old row 3 (old rows 3..10 shift to 5..12)
$ws1.Rows.Item(3).Resize(2).Insert()

# New row 5 holds what used to be row 3's data (id/name/price/existTime only)
$ws1.Range("A5").Value = 10002
$ws1.Range("B5").Value = "中国物品2"
$ws1.Range("D5").Value = 20000
$ws1.Range("F5").Value = 40

# Row 6 (old row 4 shifted down) - bornPoint column gets combined value
$ws1.Range("E6").Value = "出生点3|出生点4|出生点5"

# Fix column A width to match sheet "England" (auto best-fit, narrow)
$ws1.Columns.Item(1).ColumnWidth = 7.5

# Selection / active cell bookkeeping to match the authored edit
$ws1.Range("A3:XFD3").Select()

# ---------------------------------------------------------------------------
# Sheet "England" (sheet2) ----------------------------------------------------
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("England")

# Header E1 typo fix: bornPoint.float.array.aaa -> bornPoint.float.array
$ws2.Range("E1").Value = "bornPoint.float.array"

# Selection bookkeeping to match the authored edit
$ws2.Range("B34").Select()
